# short term societal values updated
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1) Header text tweaks: spaces -> underscores
$ws.Range("H1").Value = "selected_value"
$ws.Range("M1").Value = "source_dollar_year"

# 2) Row 7 (income_reduced / 10-year value to detainee / Reduced income):
#    Min/Mean/Median/Max were stored as text representations of negative
#    numbers ("-3,677" etc.) - replace with actual positive numeric values.
$ws.Range("H7").Value = 3677
$ws.Range("I7").Value = 1710
$ws.Range("K7").Value = 3677
$ws.Range("L7").Value = 5623

# 3) Flip sign on a few other "selected value" cells that were stored
#    as negative numbers -> now stored as positive numbers.
$ws.Range("H12").Value = 99.44
$ws.Range("H13").Value = 11
$ws.Range("H14").Value = 249634

# 4) New value for n_society (row 20): selected_value = 5,171,000,
#    formatted with thousands separators.
$ws.Range("H20").Value = 5171000
$ws.Range("H20").NumberFormat = "#,##0"

# 5) View state: zoom + selection changes
$ws.Application.ActiveWindow.Zoom = 75
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Rows.Item(1).Select()
